$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Superscript citation "1-4" -> "1,2" in the introduction paragraph.
# ------------------------------------------------------------------
$introPara = $d.Paragraphs(2)
$introStart = $introPara.Range.Start
$introText = $introPara.Range.Text
$dash = [char]0x2013
$needle = [string]::Concat("1", $dash, "4")
$idx = $introText.IndexOf($needle)
if ($idx -ge 0) {
    $target = $d.Range($introStart + $idx, $introStart + $idx + $needle.Length)
    # Superscript (w:vertAlign) survives a plain Range.Text assignment in
    # this runtime, unlike Bold/Italic, so no need to re-apply it here.
    $target.Text = "1,2"
}

# ------------------------------------------------------------------
# 2. Rewrite bibliography entry (2): Rugar et al., Nature 2004 ->
#    Kuehn et al., J. Chem. Phys. 2008, keeping the run-level
#    formatting (italic journal, bold year, italic volume, hyperlink).
# ------------------------------------------------------------------
$bib2 = $d.Paragraphs(9)
$bib2Start = $bib2.Range.Start

function Set-RelRange($paraStart, $relStart, $relEnd, $newText, $italic, $bold) {
    $rng = $d.Range($paraStart + $relStart, $paraStart + $relEnd)
    $rng.Text = $newText
    if ($italic) { $rng.Font.Italic = 1 }
    if ($bold) { $rng.Font.Bold = 1 }
}

# Original (relative) run offsets within paragraph 9:
#   0-116   "(2) Rugar, ... microscopy."        (plain)
#   117-123 "Nature"                             (italic)
#   124-128 "2004"                               (bold)
#   130-133 "430"                                (italic)
#   134-154 "(6997), 329-332 DOI:"               (plain)
# Replace from the end backwards so earlier offsets stay valid.
Set-RelRange $bib2Start 134 154 "(5), 052208 DOI:" $false $false
Set-RelRange $bib2Start 130 133 "128" $true $false
Set-RelRange $bib2Start 124 128 "2008" $false $true
Set-RelRange $bib2Start 117 123 "J. Chem. Phys." $true $false
Set-RelRange $bib2Start 0 116 "(2) Kuehn, S.; Hickman, S. A.; Marohn, J. A. Advances in mechanical detection of magnetic resonance." $false $false

# Update the DOI hyperlink text + target via the Hyperlinks collection
# (direct Range addressing does not work reliably across hyperlink field
# boundaries).
$hlink = $d.Hyperlinks.Item(1)
$hlink.TextToDisplay = "10.1063/1.2834737"
$hlink.Address = "https://doi.org/10.1063/1.2834737"

# ------------------------------------------------------------------
# 3. Remove bibliography entries (3) and (4) entirely.
# ------------------------------------------------------------------
$p3 = $d.Paragraphs(10)
$p4 = $d.Paragraphs(11)
$delRange = $d.Range($p3.Range.Start, $p4.Range.End)
$delRange.Delete()
